$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (existing rows 2-17 shift down to 3-18),
# carrying the formatting of the row being pushed down (row 2, "Isa").
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the "OpenValue" entry.
$ws.Range("A2").Value = "OpenValue"
$ws.Range("C2").Value = "8 + 16"
$ws.Range("D2").Value = 50
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Update two existing amounts that changed ("Victor" and "Bill" rows,
# now at rows 12 and 13 after the insert above).
$ws.Range("D12").Value = 39
$ws.Range("D13").Value = 30

# Update the active selection to match the new state.
$ws.Range("J14").Select()
